$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header cells, I1 ("I0") and J1 ("IF"), copying the
# formatting (bold, border, centered) already used by the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill the new columns for the data rows: I is a constant 1, J duplicates
# the existing "IP" (column H) value for that row.
$lastRow = 35
for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
